$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1099.0785
$ws.Range("J17").Value = 899.8
$ws.Range("L17").Value = 2699.4
$ws.Range("N17").Value = -3035.4

# Row 29
$ws.Range("H29").Value = 9900
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""

# Row 32
$ws.Range("H32").Value = 1001
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1001
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1001
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -1653

# Row 41
$ws.Range("H41").Value = 333.6154
$ws.Range("I41").Value = 255.25
$ws.Range("K41").Value = 255.25
$ws.Range("M41").Value = 184.75

# Row 43
$ws.Range("H43").Value = 1309.8125
$ws.Range("J43").Value = 1383.2307
$ws.Range("L43").Value = 1383.2307
$ws.Range("N43").Value = -1521.2307

# Row 51
$ws.Range("H51").Value = 4984.4
$ws.Range("J51").Value = 5110.5
$ws.Range("L51").Value = 5110.5
$ws.Range("N51").Value = -6078.5

# Row 55
$ws.Range("H55").Value = 317.1
$ws.Range("I55").Value = 263.33334
$ws.Range("J55").Value = 397.75
$ws.Range("K55").Value = 263.33334
$ws.Range("L55").Value = 397.75
$ws.Range("M55").Value = -49.33334000000002
$ws.Range("N55").Value = -825.75

# Row 70
$ws.Range("H70").Value = 50444
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 50444
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 151332
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -151872

# Row 73
$ws.Range("H73").Value = 50444
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 50444
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 151332
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -153204

# Row 82
$ws.Range("H82").Value = 3500
$ws.Range("I82").Value = 3500
$ws.Range("K82").Value = 10500
$ws.Range("M82").Value = -10094

# Row 85
$ws.Range("H85").Value = 3500
$ws.Range("I85").Value = 3500
$ws.Range("K85").Value = 10500
$ws.Range("M85").Value = -9096

# Row 132
$ws.Range("H132").Value = 984.1951
$ws.Range("I132").Value = 988.8
$ws.Range("K132").Value = 2966.4
$ws.Range("M132").Value = -436.3999999999996

# Row 137
$ws.Range("H137").Value = 1724
$ws.Range("I137").Value = 1281.3636
$ws.Range("J137").Value = 2166.6365
$ws.Range("K137").Value = 3844.0908
$ws.Range("L137").Value = 6499.9095
$ws.Range("M137").Value = -1294.0908
$ws.Range("N137").Value = -11599.9095

# Row 138
$ws.Range("H138").Value = 1754.5
$ws.Range("I138").Value = 1743.8334
$ws.Range("J138").Value = 1764.1
$ws.Range("K138").Value = 5231.5002
$ws.Range("L138").Value = 5292.299999999999
$ws.Range("M138").Value = -91.5002000000004
$ws.Range("N138").Value = -15572.3

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6259.0967
$ws.Range("I32").Value = 4467.1904
$ws.Range("K32").Value = 4467.1904
$ws.Range("M32").Value = -4180.1904

# Row 45
$ws.Range("H45").Value = 2239.7273
$ws.Range("I45").Value = 1034.2
$ws.Range("K45").Value = 1034.2
$ws.Range("M45").Value = -657.2

# Row 61
$ws.Range("H61").Value = 2230.1482
$ws.Range("I61").Value = 1324.7391
$ws.Range("K61").Value = 1324.7391
$ws.Range("M61").Value = -1112.7391

# Row 132
$ws.Range("H132").Value = 2103.8667
$ws.Range("I132").Value = 1671.6666
$ws.Range("K132").Value = 5014.9998
$ws.Range("M132").Value = -2484.9998

# Row 136
$ws.Range("H136").Value = 2230.1482
$ws.Range("I136").Value = 1324.7391
$ws.Range("K136").Value = 3974.2173
$ws.Range("M136").Value = -1424.2173

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 630.3
$ws.Range("I64").Value = 661.2
$ws.Range("J64").Value = 599.4
$ws.Range("K64").Value = 661.2
$ws.Range("L64").Value = 599.4
$ws.Range("M64").Value = -436.2
$ws.Range("N64").Value = -1049.4

# Row 67
$ws.Range("H67").Value = 630.3
$ws.Range("I67").Value = 661.2
$ws.Range("J67").Value = 599.4
$ws.Range("K67").Value = 661.2
$ws.Range("L67").Value = 599.4
$ws.Range("M67").Value = 118.8
$ws.Range("N67").Value = -2159.4

# Row 80
$ws.Range("H80").Value = 7210.9287
$ws.Range("I80").Value = 35.666668
$ws.Range("J80").Value = 9167.817999999999
$ws.Range("K80").Value = 35.666668
$ws.Range("L80").Value = 9167.817999999999
$ws.Range("M80").Value = 962.333332
$ws.Range("N80").Value = -11163.818

# Row 83
$ws.Range("H83").Value = 7210.9287
$ws.Range("I83").Value = 35.666668
$ws.Range("J83").Value = 9167.817999999999
$ws.Range("K83").Value = 178.33334
$ws.Range("L83").Value = 45839.09
$ws.Range("M83").Value = 4813.66666
$ws.Range("N83").Value = -55823.09

# Row 86
$ws.Range("H86").Value = 106753.58
$ws.Range("I86").Value = 1500.2727
$ws.Range("J86").Value = 251476.88
$ws.Range("K86").Value = 1500.2727
$ws.Range("L86").Value = 251476.88
$ws.Range("M86").Value = -377.2727
$ws.Range("N86").Value = -253722.88

# Row 89
$ws.Range("H89").Value = 106753.58
$ws.Range("I89").Value = 1500.2727
$ws.Range("J89").Value = 251476.88
$ws.Range("K89").Value = 7501.363499999999
$ws.Range("L89").Value = 1257384.4
$ws.Range("M89").Value = -1885.363499999999
$ws.Range("N89").Value = -1268616.4

# Row 94
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("K94").Value = 2000
$ws.Range("M94").Value = -1549

# Row 95
$ws.Range("H95").Value = 71896.25
$ws.Range("J95").Value = 71896.25
$ws.Range("L95").Value = 71896.25
$ws.Range("N95").Value = -77388.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2824.2917
$ws.Range("I31").Value = 1351.9412
$ws.Range("J31").Value = 6400
$ws.Range("K31").Value = 1351.9412
$ws.Range("L31").Value = 6400
$ws.Range("M31").Value = -1056.9412
$ws.Range("N31").Value = -6990

# Row 34
$ws.Range("H34").Value = 2824.2917
$ws.Range("I34").Value = 1351.9412
$ws.Range("J34").Value = 6400
$ws.Range("K34").Value = 1351.9412
$ws.Range("L34").Value = 6400
$ws.Range("M34").Value = -1149.9412
$ws.Range("N34").Value = -6804

# Row 58
$ws.Range("H58").Value = 1157.5385
$ws.Range("I58").Value = 867.05554
$ws.Range("K58").Value = 867.05554
$ws.Range("M58").Value = -664.05554

# Row 94
$ws.Range("H94").Value = 1055.2858
$ws.Range("I94").Value = 921.6667
$ws.Range("J94").Value = 1155.5
$ws.Range("K94").Value = 921.6667
$ws.Range("L94").Value = 1155.5
$ws.Range("M94").Value = -470.6667
$ws.Range("N94").Value = -2057.5

# Row 136
$ws.Range("H136").Value = 1157.5385
$ws.Range("I136").Value = 867.05554
$ws.Range("K136").Value = 2601.16662
$ws.Range("M136").Value = -51.16661999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 287
$ws.Range("I26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("M26").Value = -612

# Row 33
$ws.Range("H33").Value = 90.59999999999999
$ws.Range("I33").Value = 124.833336
$ws.Range("J33").Value = 39.25
$ws.Range("K33").Value = 749.000016
$ws.Range("L33").Value = 235.5
$ws.Range("M33").Value = -466.000016
$ws.Range("N33").Value = -801.5

# Row 107
$ws.Range("H107").Value = 743.8
$ws.Range("I107").Value = 303
$ws.Range("J107").Value = 775.2857
$ws.Range("K107").Value = 909
$ws.Range("L107").Value = 2325.8571
$ws.Range("M107").Value = 1011
$ws.Range("N107").Value = -6165.8571

# Row 131
$ws.Range("H131").Value = 793.97
$ws.Range("J131").Value = 796.9596
$ws.Range("L131").Value = 2390.8788
$ws.Range("N131").Value = -12470.8788

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2148.8572
$ws.Range("I122").Value = 1909.8572
$ws.Range("J122").Value = 2626.8572
$ws.Range("K122").Value = 5729.571599999999
$ws.Range("L122").Value = 7880.571599999999
$ws.Range("M122").Value = -3279.571599999999
$ws.Range("N122").Value = -12780.5716

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3705
$ws.Range("I16").Value = 3705
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3705
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3535
$ws.Range("N16").Value = ""

# Row 22
$ws.Range("H22").Value = 2300
$ws.Range("I22").Value = 1950
$ws.Range("K22").Value = 1950
$ws.Range("M22").Value = -1655

# Row 27
$ws.Range("H27").Value = 2300
$ws.Range("I27").Value = 1950
$ws.Range("K27").Value = 1950
$ws.Range("M27").Value = -1843

# Row 46
$ws.Range("H46").Value = 1636.4
$ws.Range("I46").Value = 1288
$ws.Range("K46").Value = 1288
$ws.Range("M46").Value = -1100

# Row 55
$ws.Range("H55").Value = 266.14285
$ws.Range("I55").Value = 281.6
$ws.Range("J55").Value = 257.55554
$ws.Range("K55").Value = 281.6
$ws.Range("L55").Value = 257.55554
$ws.Range("M55").Value = -108.6
$ws.Range("N55").Value = -603.5555400000001

# Row 95
$ws.Range("H95").Value = 80000
$ws.Range("J95").Value = 80000
$ws.Range("L95").Value = 80000
$ws.Range("N95").Value = -85492

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""

# Row 132
$ws.Range("H132").Value = 1980.6666
$ws.Range("I132").Value = 1422
$ws.Range("J132").Value = 2315.8667
$ws.Range("K132").Value = 4266
$ws.Range("L132").Value = 6947.6001
$ws.Range("M132").Value = -1736
$ws.Range("N132").Value = -12007.6001

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 27666.334
$ws.Range("J54").Value = 28999.5
$ws.Range("L54").Value = 28999.5
$ws.Range("N54").Value = -30039.5

# Row 132
$ws.Range("H132").Value = 2763.9333
$ws.Range("I132").Value = 2519.64
$ws.Range("J132").Value = 3985.4
$ws.Range("K132").Value = 7558.92
$ws.Range("L132").Value = 11956.2
$ws.Range("M132").Value = -5028.92
$ws.Range("N132").Value = -17016.2
